$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in B6:B13 from 0 to 20
$ws.Range("B6:B13").Value = 20

# Move the active selection to B14
$ws.Range("B14").Select()
